$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 118, pushing existing rows 118-145 down to 119-146
$ws.Rows("118:118").Insert()

# Populate the newly inserted row 118 with the new record
$ws.Range("A118").Value = 4
$ws.Range("B118").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C118").Value = 'Los Lagos'
$ws.Range("D118").Value = 44508
$ws.Range("E118").Value = 10
$ws.Range("F118").Value = 100112032
$ws.Range("G118").Value = 'Zapallo italiano'
$ws.Range("H118").Value = 'Sin especificar'
$ws.Range("I118").Value = 'Primera'
$ws.Range("J118").Value = 80
$ws.Range("K118").Value = 11000
$ws.Range("L118").Value = 11000
$ws.Range("M118").Value = 11000
$ws.Range("N118").Value = '$/caja 50 unidades'
$ws.Range("O118").Value = 'Región de O''Higgins'
$ws.Range("P118").Value = 220
$ws.Range("Q118").Value = 50
$ws.Range("R118").Value = 'Hortaliza'
